$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (matches the workbook's activeTab moving to this sheet)
$ws.Activate()

# Insert a new (blank) column before column N (14th column), shifting the
# existing "Late" / "heading" / "Outstanding" columns one to the right.
$ws.Columns.Item(14).Insert()

# The newly inserted column inherits the width of its left neighbour (column M).
$ws.Columns.Item(14).ColumnWidth = 9.83

# Leave the selection on cell R7, matching where editing finished.
$ws.Range("R7").Select()
